$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.598.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.320.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.31%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.82%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.318.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.22%  "

$ws.Range("E9").Value = "  +0.70%  "

$ws.Range("E10").Value = "  +2.51%  "

$ws.Range("E11").Value = "  +3.30%  "

$ws.Range("E12").Value = "  +2.06%  "

$ws.Range("E13").Value = "  +0.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.866.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.27%  "

$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.319.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.696.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.85%  "

$ws.Range("E19").Value = "  +2.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.15"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.735"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.50%  "

$ws.Range("E23").Value = "  +5.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.36%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.15%  "

$ws.Range("B28").Value = "FirstDigitalUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.80%  "

$ws.Range("E30").Value = "  +3.07%  "

$ws.Range("E31").Value = "  +2.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.08%  "

$ws.Range("E33").Value = "  +0.92%  "

$ws.Range("E34").Value = "  -0.57%  "

$ws.Range("E35").Value = "  +3.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0744"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0399"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "433.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.093.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.56%  "

$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.49%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("E44").Value = "  +4.75%  "

$ws.Range("E45").Value = "  +0.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.42%  "

$ws.Range("E49").Value = "  -0.07%  "

$ws.Range("E50").Value = "  +0.94%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.98%  "
